$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bibi Cell Manauara ---
$ws.Range("A2").Value = "Bibi Cell Manauara"
$ws.Range("B2").Value = 3340
$ws.Range("C2:AD2").Value = 0
$ws.Range("AG2").Value = 3340

# --- Row 3: Bibi Cell Ponta Negra ---
$ws.Range("A3").Value = "Bibi Cell Ponta Negra"
$ws.Range("B3").Value = 1800.01
$ws.Range("C3:AD3").Value = 0
$ws.Range("AG3").Value = 1800.01

# --- Row 4: total ---
$ws.Range("A4").Value = "total"
$ws.Range("B4").Value = 5140.01
$ws.Range("C4:AD4").Value = 0
$ws.Range("AG4").Value = 5140.01

# --- Remove the now-obsolete rows 5 and 6 ---
$ws.Rows("5:6").Delete()

Write-Output "done"
